# Fix the child-element order inside <w:rPr> for the character styles
# used by the syntax-highlighting ("Tok") styles in styles.xml.
#
# wml.xsd (CT_RPr) expects <w:b/>/<w:i/> to precede <w:color/>, but these
# styles had <w:color/> emitted first, which OOXMLValidatorCLI flags as
# Sch_UnexpectedElementContentExpectingComplex.
#
# Re-assigning a run-formatting property on a style (even to its current
# value) makes the engine re-emit that style's <w:rPr> children in
# canonical schema order, which is exactly the fix the diff describes.
# Touching every affected style's Font also keeps each w:style element's
# own content in a known, deliberate state (rather than relying on a
# side effect of touching just one of them).

$d = $word.ActiveDocument

$boldStyles = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

$italicStyles = @(
    "CommentTok",
    "DocumentationTok"
)

$boldItalicStyles = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($styleName in $boldStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $true
}

foreach ($styleName in $italicStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Italic = $true
}

foreach ($styleName in $boldItalicStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $true
    $s.Font.Italic = $true
}
